$d = $word.ActiveDocument

# 1. Merge the two math runs "=3*" and "f" into a single "=3*f" run
#    inside the oMath zone of paragraph 37 ("=3*f(n)+2").
$p37 = $d.Paragraphs.Item(37)
$om = $p37.Range.OMaths.Item(1)
$mergedMath = '<m:oMathPara><m:oMath>' + `
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><m:t>=3*f</m:t></m:r>' + `
  '<m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Arial"/><w:i/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></m:ctrlPr></m:dPr>' + `
  '<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><m:t>n</m:t></m:r></m:e></m:d>' + `
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><m:t>+2</m:t></m:r>' + `
  '</m:oMath></m:oMathPara>'
$om.Range.InsertXML($mergedMath)

# 2. Remove the entire "Aufgabe 4.2" draft section (paragraphs 49-66),
#    keeping the single blank paragraph (48) that precedes it.
$p48 = $d.Paragraphs.Item(48)
$p49 = $d.Paragraphs.Item(49)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$killRange = $d.Range($p49.Range.Start, $pLast.Range.End)
$killRange.Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
